# Updates to MC and alignment
# Add a new "F" / 500 column pair (O, P) to every data row, update the
# current selection, and set the page setup for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 40

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "F"
    $ws.Cells.Item($r, 16).Value = 500
}

# Update page setup (paper size A4, portrait orientation)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Update the active selection
$ws.Range("N33:N40").Select() | Out-Null
